# Refresh cached Universalis market-board pull for the Leve profit tables.
# Columns: H currentAveragePrice | I currentAveragePriceNQ | J currentAveragePriceHQ
#          K LevePriceNQ | L LevePriceHQ | M LeveProfitNQ | N LeveProfitHQ
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

# Row 17: One for the Road / Potion
$ws.Range("H17").Value = 1937.2273
$ws.Range("J17").Value = 1937.2273
$ws.Range("L17").Value = 5811.6819
$ws.Range("N17").Value = -6147.6819

# Row 70: Consecrating Congregation / Holy Water
$ws.Range("H70").Value = 4376.8184
$ws.Range("J70").Value = 4093.125
$ws.Range("L70").Value = 12279.375
$ws.Range("N70").Value = -12819.375

# Row 73: Curbing the Contagion (L) / Holy Water
$ws.Range("H73").Value = 4376.8184
$ws.Range("J73").Value = 4093.125
$ws.Range("L73").Value = 12279.375
$ws.Range("N73").Value = -14151.375

# Row 100: Asking for a Friend / Beetle Glue
$ws.Range("H100").Value = 2500
$ws.Range("I100").Value = 2500
$ws.Range("J100").Value = 2500
$ws.Range("K100").Value = 2500
$ws.Range("L100").Value = 2500
$ws.Range("M100").Value = -1959
$ws.Range("N100").Value = -3582

# Row 106: Making Your Mark / Enchanted Palladium Ink
$ws.Range("H106").Value = 66669668
$ws.Range("I106").Value = 111112450
$ws.Range("K106").Value = 111112450
$ws.Range("M106").Value = -111111819

# Row 112: Making Ends Meet / Superior Spiritbond Potion
$ws.Range("H112").Value = 1521
$ws.Range("J112").Value = 1543.3448
$ws.Range("L112").Value = 4630.0344
$ws.Range("N112").Value = -6846.0344

# Row 113: Amaro Kart / Starch Glue
$ws.Range("H113").Value = 3511.1428
$ws.Range("I113").Value = 2835.7144
$ws.Range("J113").Value = 3848.8572
$ws.Range("K113").Value = 2835.7144
$ws.Range("L113").Value = 3848.8572
$ws.Range("M113").Value = 418.2856000000002
$ws.Range("N113").Value = -10356.8572

# Row 116: Growing Up / Growth Formula Kappa
$ws.Range("H116").Value = 34728200
$ws.Range("J116").Value = 125005750
$ws.Range("L116").Value = 125005750
$ws.Range("N116").Value = -125012634


$ws = $wb.Worksheets.Item("ARM")

# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 13034.973
$ws.Range("I32").Value = 13759.5
$ws.Range("J32").Value = 10499.125
$ws.Range("K32").Value = 13759.5
$ws.Range("L32").Value = 10499.125
$ws.Range("M32").Value = -13472.5
$ws.Range("N32").Value = -11073.125

# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 7354226.5
$ws.Range("J74").Value = 1689.8572
$ws.Range("L74").Value = 1689.8572
$ws.Range("N74").Value = -3437.8572

# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 7354226.5
$ws.Range("J77").Value = 1689.8572
$ws.Range("L77").Value = 8449.286
$ws.Range("N77").Value = -17185.286

# Row 86: Sir, Dost Thou Even Heft / Adamantite Chain Hose of Fending
$ws.Range("H86").Value = 30000
$ws.Range("J86").Value = 30000
$ws.Range("L86").Value = 30000
$ws.Range("N86").Value = -32372

# Row 89: Men in Adamantite (L) / Adamantite Chain Hose of Fending
$ws.Range("H89").Value = 30000
$ws.Range("J89").Value = 30000
$ws.Range("L89").Value = 90000
$ws.Range("N89").Value = -101856

# Row 122: Haste for High Durium / High Durium Nugget
$ws.Range("H122").Value = 4871.857
$ws.Range("I122").Value = 4871.857
$ws.Range("K122").Value = 14615.571
$ws.Range("M122").Value = -12165.571

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 4966.0464
$ws.Range("I132").Value = 2387.15
$ws.Range("K132").Value = 7161.450000000001
$ws.Range("M132").Value = -4631.450000000001


$ws = $wb.Worksheets.Item("BSM")

# Row 76: Keep Up with the Mechanics / Titanium-barreled Arquebus
$ws.Range("H76").Value = 66607
$ws.Range("J76").Value = 66607
$ws.Range("L76").Value = 66607
$ws.Range("N76").Value = -67237

# Row 79: Unconventional Weaponry (L) / Titanium-barreled Arquebus
$ws.Range("H79").Value = 66607
$ws.Range("J79").Value = 66607
$ws.Range("L79").Value = 66607
$ws.Range("N79").Value = -68791

# Row 86: Through Thick and Thin / Adamantite Nugget
$ws.Range("H86").Value = 2149.1765
$ws.Range("I86").Value = 1987.8572
$ws.Range("K86").Value = 1987.8572
$ws.Range("M86").Value = -864.8571999999999

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget
$ws.Range("H89").Value = 2149.1765
$ws.Range("I89").Value = 1987.8572
$ws.Range("K89").Value = 9939.286
$ws.Range("M89").Value = -4323.286

# Row 105: Ingot to Wing It / Molybdenum Ingot
$ws.Range("H105").Value = 2245.625
$ws.Range("I105").Value = 2027.579
$ws.Range("K105").Value = 2027.579
$ws.Range("M105").Value = -280.579


$ws = $wb.Worksheets.Item("CRP")

# Row 16: Raise the Roof / Ash Lumber
$ws.Range("H16").Value = 1926.6923
$ws.Range("I16").Value = 1777
$ws.Range("K16").Value = 1777
$ws.Range("M16").Value = -1490

# Row 22: Driving Up the Wall / Elm Lumber
$ws.Range("H22").Value = 912
$ws.Range("I22").Value = 475.66666
$ws.Range("J22").Value = 1348.3334
$ws.Range("K22").Value = 475.66666
$ws.Range("L22").Value = 1348.3334
$ws.Range("M22").Value = -125.66666
$ws.Range("N22").Value = -2048.3334

# Row 58: You Do the Heavy Lifting / Mahogany Lumber
$ws.Range("H58").Value = 4139.3335
$ws.Range("J58").Value = 7320.357
$ws.Range("L58").Value = 7320.357
$ws.Range("N58").Value = -7726.357

# Row 94: Beech, Please / Beech Lumber
$ws.Range("H94").Value = 4982
$ws.Range("J94").Value = 4982
$ws.Range("L94").Value = 4982
$ws.Range("N94").Value = -5884

# Row 107: Built to Last / White Oak Lumber
$ws.Range("H107").Value = 3636964.5
$ws.Range("I107").Value = 9091154
$ws.Range("J107").Value = 837.6667
$ws.Range("K107").Value = 9091154
$ws.Range("L107").Value = 837.6667
$ws.Range("M107").Value = -9089234
$ws.Range("N107").Value = -4677.6667

# Row 113: Patient Patients / White Ash Lumber
$ws.Range("H113").Value = 1926.6923
$ws.Range("I113").Value = 1777
$ws.Range("K113").Value = 1777
$ws.Range("M113").Value = 393

# Row 136: Turali Quality / Dark Mahogany Lumber
$ws.Range("H136").Value = 4139.3335
$ws.Range("J136").Value = 7320.357
$ws.Range("L136").Value = 21961.071
$ws.Range("N136").Value = -27061.071


$ws = $wb.Worksheets.Item("CUL")

# Row 123: Topping Up the Pot / Zurek
$ws.Range("H123").Value = 19250
$ws.Range("I123").Value = 13500
$ws.Range("K123").Value = 40500
$ws.Range("M123").Value = -38050

# Row 134: Don't Knock It Till You've Tried It / Mezcal-marinated Swampmonk
$ws.Range("H134").Value = 10483.583
$ws.Range("I134").Value = 3336.2222
$ws.Range("J134").Value = 31925.666
$ws.Range("K134").Value = 10008.6666
$ws.Range("L134").Value = 95776.99800000001
$ws.Range("M134").Value = -4938.6666
$ws.Range("N134").Value = -105916.998

# Row 139: Najoothie / Wild Banana Blend
$ws.Range("H139").Value = 2623.9546
$ws.Range("I139").Value = 1266.5
$ws.Range("J139").Value = 4999.5
$ws.Range("K139").Value = 3799.5
$ws.Range("L139").Value = 14998.5
$ws.Range("M139").Value = 1340.5
$ws.Range("N139").Value = -25278.5


$ws = $wb.Worksheets.Item("GSM")

# Row 102: Put the Metal to the Peddle / Durium Ingot
$ws.Range("H102").Value = 26325600
$ws.Range("I102").Value = 35724172
$ws.Range("K102").Value = 35724172
$ws.Range("M102").Value = -35722550

# Row 122: Awarding Academic Excellence / Ametrine
$ws.Range("H122").Value = 444935.62
$ws.Range("I122").Value = 788028.9399999999
$ws.Range("K122").Value = 2364086.82
$ws.Range("M122").Value = -2361636.82

# Row 126: Gold Rush Order / Phrygian Gold Ingot
$ws.Range("H126").Value = 8999.75
$ws.Range("I126").Value = 2999.5
$ws.Range("K126").Value = 8998.5
$ws.Range("M126").Value = -6528.5


$ws = $wb.Worksheets.Item("LTW")

# Row 22: Skin off Their Backs / Aldgoat Leather
$ws.Range("H22").Value = 933.13336
$ws.Range("I22").Value = 783
$ws.Range("J22").Value = 1158.3334
$ws.Range("K22").Value = 783
$ws.Range("L22").Value = 1158.3334
$ws.Range("M22").Value = -488
$ws.Range("N22").Value = -1748.3334

# Row 27: Fire and Hide / Aldgoat Leather
$ws.Range("H27").Value = 933.13336
$ws.Range("I27").Value = 783
$ws.Range("J27").Value = 1158.3334
$ws.Range("K27").Value = 783
$ws.Range("L27").Value = 1158.3334
$ws.Range("M27").Value = -676
$ws.Range("N27").Value = -1372.3334

# Row 46: Supply Side Logic / Boar Leather
$ws.Range("H46").Value = 6844.391
$ws.Range("J46").Value = 6991
$ws.Range("L46").Value = 6991
$ws.Range("N46").Value = -7367

# Row 55: It's Not a Job, It's a Calling / Peiste Leather
$ws.Range("H55").Value = 330.8889
$ws.Range("I55").Value = 141.66667
$ws.Range("K55").Value = 141.66667
$ws.Range("M55").Value = 31.33332999999999

# Row 93: Hide to Go Seek / Gagana Leather
$ws.Range("H93").Value = 2962.3
$ws.Range("J93").Value = 1127.6666
$ws.Range("L93").Value = 1127.6666
$ws.Range("N93").Value = -3623.6666

# Row 117: If I Could Walk a Thousand Malms / Zonureskin Shoes of Gathering
$ws.Range("H117").Value = 49500
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

# Row 122: Hell on Leather / Gaja Leather
$ws.Range("H122").Value = 62506056
$ws.Range("I122").Value = 125004360
$ws.Range("K122").Value = 375013080
$ws.Range("M122").Value = -375010630

# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 4597.6562
$ws.Range("I132").Value = 3536.55
$ws.Range("K132").Value = 10609.65
$ws.Range("M132").Value = -8079.650000000001


$ws = $wb.Worksheets.Item("WVR")

# Row 2: The Unmentionables / Hempen Underpants
$ws.Range("H2").Value = 8031.4165
$ws.Range("I2").Value = 8031.4165
$ws.Range("K2").Value = 8031.4165
$ws.Range("M2").Value = -7919.4165

# Row 4: Not Cool Enough / Hempen Undershirt
$ws.Range("H4").Value = 91.5
$ws.Range("I4").Value = 92.22221999999999
$ws.Range("J4").Value = 85
$ws.Range("K4").Value = 92.22221999999999
$ws.Range("L4").Value = 85
$ws.Range("M4").Value = 20.77778000000001
$ws.Range("N4").Value = -311

# Row 27: Hitting Below the Belt / Cotton Breeches of Crafting
$ws.Range("H27").Value = 75626
$ws.Range("J27").Value = 75626
$ws.Range("L27").Value = 75626
$ws.Range("N27").Value = -75764

# Row 37: Bet You Anything / Velveteen Sarouel of Gathering
$ws.Range("H37").Value = 39994.5
$ws.Range("J37").Value = 49999
$ws.Range("L37").Value = 49999
$ws.Range("N37").Value = -50405

# Row 45: Private Concerns / Linen Trousers
$ws.Range("H45").Value = 69600
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 69600
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 69600
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -70582

# Row 100: Of Great Import / Kudzu Thread
$ws.Range("H100").Value = 1016606.2
$ws.Range("J100").Value = 4166.5
$ws.Range("L100").Value = 8333
$ws.Range("N100").Value = -9415

# Row 102: Don't Sweat the Role / Serge Turban of Crafting
$ws.Range("H102").Value = 100000
$ws.Range("J102").Value = 100000
$ws.Range("L102").Value = 100000
$ws.Range("N102").Value = -106490

# Row 113: A Tender Table / Pixie Floss
$ws.Range("H113").Value = 973.5625
$ws.Range("I113").Value = 969.7
$ws.Range("J113").Value = 980
$ws.Range("K113").Value = 2909.1
$ws.Range("L113").Value = 2940
$ws.Range("M113").Value = -739.1000000000004
$ws.Range("N113").Value = -7280

# Row 127: Turban Sprawl / Snow Linen Turban of Crafting
$ws.Range("H127").Value = 58700
$ws.Range("J127").Value = 58700
$ws.Range("L127").Value = 58700
$ws.Range("N127").Value = -68620

# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 10873040
$ws.Range("I132").Value = 1217.6786
$ws.Range("K132").Value = 3653.0358
$ws.Range("M132").Value = -1123.0358

